# Finalized definition of virtual environment
# Updates calibration parameter values (and their abs_epsi_autocorr arrays)
# across the AR, SETAR, GARCH, TARCH and AR_TARCH sheets.

$wb = $excel.ActiveWorkbook

# --- AR sheet ---
$wsAR = $wb.Worksheets.Item("AR")
$wsAR.Range("B2").Value = 0.008916311388064143
$wsAR.Range("B3").Value = 0.8328228104379528
$wsAR.Range("B4").Value = 67.82037903986348
$wsAR.Range("B5").Value = "[0.9999999999999999, 0.2797259506740059, 0.24900086961075893, 0.26949077737715676, 0.2618909636715644, 0.426479110316983, 0.2665015360641186, 0.2399550456955552, 0.22145681171563353, 0.22916581144241424, 0.1845653397671616, 0.20543203950323502, 0.18367172091167397, 0.21050976446574993, 0.22220306235239745, 0.18180292025591968, 0.18463158252494655, 0.17693472083943734, 0.16194323532437113, 0.16309724558913208]"

# --- SETAR sheet ---
$wsSETAR = $wb.Worksheets.Item("SETAR")
$wsSETAR.Range("B3").Value = 0.855294681360805
$wsSETAR.Range("B4").Value = -2.186409598444523
$wsSETAR.Range("B5").Value = 0.7959662244860162
$wsSETAR.Range("B6").Value = 55.76297824117911
$wsSETAR.Range("B7").Value = 2.462913177027644
$wsSETAR.Range("B8").Value = 0.734395547537864
$wsSETAR.Range("B9").Value = 44.32505331810656
$wsSETAR.Range("B10").Value = "[0.9999999999999998, 0.28302540337385607, 0.24129123051551393, 0.2580413616587471, 0.28938846791947126, 0.3454669991273173, 0.2385277310581725, 0.21952658077006582, 0.2001513700262536, 0.2379972378118566, 0.19025281021938578, 0.19512638235546298, 0.1839774940929245, 0.1862512062393665, 0.18120219439184418, 0.15355485208200761, 0.16721335281011498, 0.17826433311355802, 0.16017036243377217, 0.15292910102105775]"

# --- GARCH sheet ---
$wsGARCH = $wb.Worksheets.Item("GARCH")
$wsGARCH.Range("B2").Value = 0.02988465038507307
$wsGARCH.Range("B3").Value = 0.244286917739708
$wsGARCH.Range("B4").Value = 0.09536980095588932
$wsGARCH.Range("B5").Value = 0.9046301990441107
$wsGARCH.Range("B6").Value = "[1.0, -0.007575341520823182, -0.019325526781894004, -0.01681443787968801, -0.0366506087800234, 0.29539169256332315, -0.002825172575298081, 0.025208367716528058, 0.0025409853046714954, -0.0029422469368951287, -0.07616189181066904, -0.042419428338523885, -0.04010566545216876, -0.0039219378641295434, 0.03203312699760767, -0.035322937253569796, -0.000986140551119544, -0.04316362904677151, -0.04289954659174672, -0.039593526413248266]"

# --- TARCH sheet ---
$wsTARCH = $wb.Worksheets.Item("TARCH")
$wsTARCH.Range("B2").Value = -0.002319332802289378
$wsTARCH.Range("B3").Value = 0.2480089963239911
$wsTARCH.Range("B4").Value = 0.08236600999930936
$wsTARCH.Range("B5").Value = 0.9049871957104938
$wsTARCH.Range("B6").Value = "[1.0, -0.007072368391461506, -0.01886108165203145, -0.017316884045910263, -0.03723533905847705, 0.2946173805377546, -0.00206227953597705, 0.02506684424760601, 0.0028190125079451955, -0.0025894106223187576, -0.07672885664341558, -0.042048094749013255, -0.039638980977046503, -0.003520674284241749, 0.033101780356912795, -0.03506907296468009, -0.001063842940547818, -0.04263048389440927, -0.04244179584597065, -0.03928572554808717]"
$wsTARCH.Range("B7").Value = 0.02529384690701313

# --- AR_TARCH sheet ---
$wsARTARCH = $wb.Worksheets.Item("AR_TARCH")
$wsARTARCH.Range("B2").Value = 0.07554782284192292
$wsARTARCH.Range("B3").Value = 0.2093343597947837
$wsARTARCH.Range("B4").Value = 0.09283173354892066
$wsARTARCH.Range("B5").Value = 0.9100752859855565
$wsARTARCH.Range("B6").Value = "[1.0, 0.009712893278950526, -0.01693861572256187, 0.01031357864962404, -0.020627588515426735, 0.22418305593293247, -0.0060780062646489455, 0.026661235035804774, -0.020432479636356, -0.005480270973999062, -0.05683469109650953, -0.031117653915863144, -0.036466222794847, -0.00666065266668414, 0.04671001472156489, -0.04279537767878011, 0.0007994801476868401, -0.023810269584477538, -0.03615257617490866, -0.031384973727866555]"
$wsARTARCH.Range("B7").Value = -0.005814042285948937
$wsARTARCH.Range("B9").Value = 0.8154383811893925
